# Insert one new weekly record for "Vega Monumental Concepción - Acelga".
# The new record is inserted as row 69, pushing the existing rows 69:189
# down to 70:190 (dimension grows from A1:R189 to A1:R190).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 69; Excel shifts rows 69:189 -> 70:190
# and copies formatting (including the date number-format on column D)
# from the surrounding rows automatically.
$ws.Rows("69:69").Insert()

# Populate the newly inserted row 69 with the new weekly price record.
$ws.Range("A69").Value = 11
$ws.Range("B69").Value = "Vega Monumental Concepción"
$ws.Range("C69").Value = "Bíobío"
$ws.Range("D69").Value = 44546
$ws.Range("E69").Value = 8
$ws.Range("F69").Value = 100112009
$ws.Range("G69").Value = "Acelga"
$ws.Range("H69").Value = "Sin especificar"
$ws.Range("I69").Value = "Primera"
$ws.Range("J69").Value = 450
$ws.Range("K69").Value = 600
$ws.Range("L69").Value = 650
$ws.Range("M69").Value = 628
$ws.Range("N69").Value = "$/atado 0,5 a 1 kilo"
$ws.Range("O69").Value = "Región de Ñuble"
$ws.Range("P69").Value = 628
$ws.Range("Q69").Value = 1
$ws.Range("R69").Value = "Hortaliza"
